# Fix door directions on the ClueBoard sheet.
# Several door cells used a generic "<Letter>D" code; this corrects them to
# the proper directional code (matching the commit message
# "fixed door directions in xcel sheet").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H10").Value = "QR"
$ws.Range("S10").Value = "LL"
$ws.Range("S11").Value = "LL"
$ws.Range("S12").Value = "LL"
$ws.Range("D16").Value = "DL"
$ws.Range("M16").Value = "TL"
$ws.Range("D17").Value = "DL"
$ws.Range("M17").Value = "TL"
$ws.Range("Q20").Value = "GU"
$ws.Range("C21").Value = "BU"
$ws.Range("J21").Value = "RU"

# Update the saved selection to match the author's final cursor position.
$ws.Activate()
$ws.Range("D17").Select() | Out-Null
